$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "0.9990")
# are preserved verbatim as strings instead of being parsed into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.521.97"
$ws.Range("E2").Value = "  +2.13%  "
$ws.Range("D3").Value = "1.857.35"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "245.25"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "0.6967"
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("D7").Value = "0.9995"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.3076"
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07699"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").Value = "23.68"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").Value = "0.07793"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "5.171"
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("D13").Value = "1.857.84"
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.6947"
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "91.26"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").Value = "6.326"
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("D17").Value = "29.499.87"
$ws.Range("E17").Value = "  +2.09%  "
$ws.Range("D18").Value = "0.000008333"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").Value = "2.103.92"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").Value = "238.94"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").Value = "12.76"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "7.624"
$ws.Range("E23").Value = "  +2.06%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "0.1497"
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("D26").Value = "160.01"
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("D27").Value = "8.894"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").Value = "18.28"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").Value = "1.533"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("D30").Value = "4.252"
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("D31").Value = "4.153"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "1.203"
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").Value = "0.05108"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "0.7756"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").Value = "1.885"
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "1.315.22"
$ws.Range("E38").Value = "  +8.05%  "
$ws.Range("D39").Value = "0.01879"
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("D40").Value = "2.723"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").Value = "0.9522"
$ws.Range("E41").Value = "  +1.21%  "
$ws.Range("D42").Value = "106.25"
$ws.Range("E42").Value = "  -2.22%  "
$ws.Range("D43").Value = "5.779"
$ws.Range("E43").Value = "  +1.76%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").Value = "9.850"
$ws.Range("E45").Value = "  +3.39%  "
$ws.Range("E46").Value = "  +2.11%  "
$ws.Range("D47").Value = "2.001.07"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("D48").Value = "0.5233"
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("D49").Value = "1.790"
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").Value = "63.27"
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("D51").Value = "6.964"
$ws.Range("E51").Value = "  +0.97%  "

# Restore the default style so the cell style attribute matches the original
# (only the text content should differ from the source workbook).
$ws.Range("D2:D51").Style = "Normal"
